$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revert back: remove the "Buenos Aires Innovation Park (City of Buenos Aires
# Government)" row (row 4). Deleting the entire row shifts all subsequent
# rows up by one, which also removes the now-unused shared string and
# shrinks the used range from A1:AG41 to A1:AG40.
$ws.Rows.Item(4).EntireRow.Delete()
